$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the price series. It belongs right above
# the current row 137, so insert a fresh row there (this pushes the old
# rows 137-263 down to 138-264, which is exactly what the diff shows).
$ws.Rows("137").Insert()

# Populate the newly inserted row with the new observation (same static
# metadata as the surrounding rows for this market/product).
$ws.Range("A137").Value2 = 10
$ws.Range("B137").Value2 = "Vega Modelo de Temuco"
$ws.Range("C137").Value2 = "La Araucanía"
$ws.Range("D137").Value2 = 44781
$ws.Range("E137").Value2 = 9
$ws.Range("F137").Value2 = 100112043
$ws.Range("G137").Value2 = "Pepino dulce"
$ws.Range("H137").Value2 = "Cultivar IV Región"
$ws.Range("I137").Value2 = "Primera"
$ws.Range("J137").Value2 = 500
$ws.Range("K137").Value2 = 18000
$ws.Range("L137").Value2 = 19000
$ws.Range("M137").Value2 = 18600
$ws.Range("N137").Value2 = "$/bandeja 18 kilos"
$ws.Range("O137").Value2 = "Provincia de Limarí"
$ws.Range("P137").Value2 = 1033
$ws.Range("Q137").Value2 = 18
$ws.Range("R137").Value2 = "Hortaliza"

# Keep the date column's date/time number format on the new row, matching
# every other row's column D.
$ws.Range("D137").NumberFormat = $ws.Range("D138").NumberFormat
